# account_bank_statement_import_adyen / adyen_test.xlsx migration edit
# Mirrors the OOXML diff produced when this fixture was re-saved during the
# 12.0 migration: a currency value tweak, a lower-cased datetime format
# code, an updated "Gross Credit" sample value, tidied-up column widths and
# a refreshed cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "EUR" sample currency values become "USD" everywhere on the sheet ---
[void]$ws.Cells.Replace("EUR", "USD")

# --- datetime number format: YYYY-MM-DD HH:MM:SS -> yyyy-mm-dd hh:mm:ss ---
$ws.Range("G5:G25").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- sample data fix: Gross Credit (GC) on row 10 ---
$ws.Range("M10").Value = 1598

# --- column width cleanup (values rounded to match the re-saved sheet) ---
$ws.Columns.Item(3).ColumnWidth = 37.23
$ws.Columns.Item(4).ColumnWidth = 22.36
$ws.Columns.Item(5).ColumnWidth = 25.7
$ws.Columns.Item(6).ColumnWidth = 18.89
$ws.Columns.Item(7).ColumnWidth = 26.12
$ws.Columns.Item(9).ColumnWidth = 42.93
$ws.Columns.Item(10).ColumnWidth = 23.61
$ws.Columns.Item(20).ColumnWidth = 17.47
$ws.Columns.Item(22).ColumnWidth = 25.14

# --- refreshed selection/active cell ---
[void]$ws.Range("L9").Select()
